# MHD2-83: Adding requester code to report footers
# The footer currently contains the requester-code placeholder split across
# two separate runs ("REQUESTER_CODE_I" + "N"). Re-write it as a single
# contiguous run "REQUESTER_CODE_IN" via Find & Replace scoped to the
# footer range (Document.Content.Find does not reach headers/footers).

$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Sections.Count; $i++) {
    $section = $d.Sections.Item($i)
    $footers = $section.Footers
    for ($j = 1; $j -le $footers.Count; $j++) {
        $footer = $footers.Item($j)
        if ($footer.Exists) {
            $rng = $footer.Range
            $rng.Find.Execute("REQUESTER_CODE_IN", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "REQUESTER_CODE_IN", 2)
        }
    }
}
